# Auto-generated edit script applying numeric/text cell updates per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 429
$ws.Range("F5").Value = 1788
$ws.Range("F7").Value = 2878
$ws.Range("F8").Value = 2405
$ws.Range("F9").Value = 752
$ws.Range("F10").Value = 7201
$ws.Range("F11").Value = 164
$ws.Range("F13").Value = 175
$ws.Range("F14").Value = 1637
$ws.Range("F15").Value = 1404
$ws.Range("F17").Value = 124
$ws.Range("F18").Value = 124
$ws.Range("F19").Value = 3193
$ws.Range("F20").Value = 5722
$ws.Range("F21").Value = 5722
$ws.Range("F22").Value = 548
$ws.Range("F23").Value = 900
$ws.Range("F24").Value = 1168
$ws.Range("F25").Value = 320
$ws.Range("F26").Value = 5660
$ws.Range("F27").Value = 316
$ws.Range("F29").Value = 3923
$ws.Range("F31").Value = 655
$ws.Range("F32").Value = 1793
$ws.Range("F33").Value = 1112
$ws.Range("F34").Value = 232
$ws.Range("F35").Value = 11
$ws.Range("F36").Value = 125
$ws.Range("F37").Value = 50
$ws.Range("F38").Value = 297
$ws.Range("F39").Value = 1094
$ws.Range("F40").Value = 458
$ws.Range("F41").Value = 1799
$ws.Range("F42").Value = 67
$ws.Range("F43").Value = 324
$ws.Range("F44").Value = 123
$ws.Range("F45").Value = 976
$ws.Range("F47").Value = 53
$ws.Range("F48").Value = 20
$ws.Range("F49").Value = 70
$ws.Range("F50").Value = 122
$ws.Range("F51").Value = 4

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 110
$ws.Range("F8").Value = 154
$ws.Range("F10").Value = 973
$ws.Range("F14").Value = 10
$ws.Range("F16").Value = 637
$ws.Range("F17").Value = 304
$ws.Range("F21").Value = 93
$ws.Range("F23").Value = 3
$ws.Range("F25").Value = 114
$ws.Range("F28").Value = 221
$ws.Range("F33").Value = 239

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3318
$ws.Range("F7").Value = 1508
$ws.Range("F9").Value = 437
$ws.Range("F10").Value = 2963
$ws.Range("F11").Value = 376
$ws.Range("F12").Value = 708
$ws.Range("G13").Value = "不可售"
$ws.Range("F14").Value = 868
$ws.Range("F15").Value = 1378

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1508
$ws.Range("F4").Value = 429
$ws.Range("F5").Value = 437
$ws.Range("F6").Value = 2963
$ws.Range("F7").Value = 1788
$ws.Range("F8").Value = 2878
$ws.Range("F9").Value = 376
$ws.Range("F10").Value = 2405
$ws.Range("F11").Value = 752
$ws.Range("F12").Value = 7201
$ws.Range("F13").Value = 164
$ws.Range("F14").Value = 708
$ws.Range("C15").Value = "上海·第五人格0nly 2.0"
$ws.Range("D15").Value = "吴中路1588号上海爱琴海购物中心F4 竞梦元宇宙"
$ws.Range("E15").Value = "2024.07.27 10:00-07.27 17:00"
$ws.Range("F15").Value = 1637
$ws.Range("G15").Value = 68
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=88872"
$ws.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202407/5EVtnVTZ1720436658685.jpeg"
$ws.Range("C16").Value = "上海·第十二届Redamancy动漫游戏嘉年华"
$ws.Range("D16").Value = "中山北路3300号4楼L4001号 环球港上海世嘉都市乐园"
$ws.Range("E16").Value = "2024.07.27 10:00-07.28 17:00"
$ws.Range("F16").Value = 1404
$ws.Range("G16").Value = 60
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=87423"
$ws.Range("I16").Value = "//i0.hdslb.com/bfs/openplatform/202406/Ll5fKZSD1718292388904.png"
$ws.Range("C17").Value = "上海·［呪術廻戦カフェ2024第二季］主题咖啡厅"
$ws.Range("D17").Value = "西藏北路166号（地铁8号线曲阜路下） 静安大悦城"
$ws.Range("E17").Value = "2024.07.27 00:00-08.31 23:59"
$ws.Range("F17").Value = 868
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=89361"
$ws.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202407/cPSEorSS1721121089976.png"
$ws.Range("B18").Value = "2024-07-28"
$ws.Range("C18").Value = "上海·火影忍者0nly"
$ws.Range("D18").Value = "吴中路1588号上海爱琴海购物中心F4 竞梦元宇宙"
$ws.Range("E18").Value = "2024.07.28 10:00-07.28 17:00"
$ws.Range("F18").Value = 1254
$ws.Range("G18").Value = 68
$ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=88871"
$ws.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202407/EeZygAsN1720438160543.jpeg"
$ws.Range("C19").Value = "上海·音阅派国漫演唱会-《狐妖小红娘》《一人之下》领衔国漫原声音乐现场"
$ws.Range("D19").Value = "丁香路425号 上海东方艺术中心"
$ws.Range("E19").Value = "2024.07.28 19:30-07.28 21:00"
$ws.Range("F19").Value = 637
$ws.Range("G19").Value = 280
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=87560"
$ws.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202406/QuucsAfP1718693778298.jpeg"
$ws.Range("F20").Value = 124
$ws.Range("F21").Value = 124
$ws.Range("F22").Value = 1378
$ws.Range("F23").Value = 3194
$ws.Range("F24").Value = 304
$ws.Range("F25").Value = 5722
$ws.Range("F27").Value = 548
$ws.Range("F28").Value = 1168
$ws.Range("F29").Value = 320
$ws.Range("F30").Value = 5660
$ws.Range("F31").Value = 316
$ws.Range("F32").Value = 3923
$ws.Range("F33").Value = 655
$ws.Range("F35").Value = 1793
$ws.Range("F36").Value = 1112
$ws.Range("F37").Value = 232
$ws.Range("F38").Value = 11
$ws.Range("F39").Value = 114
$ws.Range("F40").Value = 125
$ws.Range("F41").Value = 50
$ws.Range("F42").Value = 1094
$ws.Range("F43").Value = 458
$ws.Range("F44").Value = 1799
$ws.Range("F45").Value = 324
$ws.Range("F46").Value = 123
$ws.Range("F47").Value = 976
$ws.Range("F49").Value = 239
$ws.Range("F50").Value = 122
$ws.Range("F51").Value = 4

